$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry added to the "journal de travail" on row 20:
# Semaine / Tache / Duree / Commentaire for the "creation de nom d'utilisateur sans fichier" task.
$ws.Range("B20").Value = 6
$ws.Range("D20").Value = "2h15"
$ws.Range("C20").Value = "Visualition de vidéos pour comment créer des fichiers dans le jeu, et création d'utilisateur avant le début de la partie"
$ws.Range("E20").Value = "J'ai visualisé des vidéos pour comment faire des fichiers et comment y écrire à l'intérieur, et le programme demande si l'utilisateur veut entrer un nom ou pas et l'affiche pendant le jeu et la fin du jeu dans le message de victoire"

# Tache cell now wraps like the Commentaire column.
$ws.Range("C20").WrapText = $true

# Row grew taller to fit the new wrapped text.
$ws.Rows.Item(20).RowHeight = 75

# A handful of earlier wrapped-text rows reflow slightly (different Excel build / fonts).
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(8).RowHeight = 60
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 45
$ws.Rows.Item(18).RowHeight = 45
$ws.Rows.Item(19).RowHeight = 165

# Selection moved on to the next empty row.
[void]$ws.Range("E21").Select()
